$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C19").Value = "[name=`"Plastic Knight`"]  ...Are you threatening me? Threatening 'Plastic' Knight Szewczyk to take on your dirty work? `n"
$ws.Range("C40").Value = "[name=`"Bald Marcin`"]  How’d 'Plastic' Szewczyk find his way to this little place? Tell you the truth, this bar doesn’t see all that many active knights in it. `n"
$ws.Range("C42").Value = "[name=`"Plastic Knight`"]  Hmph... tremendous all around. I’ll have a 'Red Edelweiss.'`n"
$ws.Range("C54").Value = "[name=`"Plastic Knight`"]  I know what you want to say, Marcin. Up there, you’re thinking, 'how’s a tiny little Plastic Knight get into such a big mess?'`n"
$ws.Range("C68").Value = "[name=`"'Flametail' Sonna`"]  Hey, you’ve been tailing us for ages now. Time you showed your face, right?`n"
$ws.Range("C69").Value = "[name=`"'Flametail' Sonna`"]  Oh, the overbearing type... No, I do NOT like fans like you...`n"
$ws.Range("C70").Value = "[name=`"'Flametail' Sonna`"]  'Wait—that bow—are you—?!' `n"
$ws.Range("C71").Value = "[name=`"'Flametail' Sonna`"]  —Thought I’d say that and keel over, Armorless Union assassin?`n"
$ws.Range("C72").Value = "[name=`"'Flametail' Sonna`"]  Hmm... so this is how you’re molded. Way you dodge about everywhere, no-one’s actually seen your looks before. Real gloomy line of work you’re in...`n"
$ws.Range("C73").Value = "[name=`"'Flametail' Sonna`"]  Don’t get the chance everyday—`n"
$ws.Range("C75").Value = "[name=`"'Flametail' Sonna`"]  —to have a real scrap, you and us.`n"
$ws.Range("C102").Value = "[name=`"Młynar`"]  And what happened? Back against the wall, she still chose the path of competition, to flaunt her 'convictions.' Yes, indeed. 'Convictions.' `n"
$ws.Range("C107").Value = "[name=`"Młynar`"]  'The Radiant Knight' is a title she should’ve never been conferred! She has no power to shoulder any glory of the sort!`n"
